$wb = $excel.ActiveWorkbook

# --- "loads" sheet (sheet3): extend header + data row with new columns ---
$wsLoads = $wb.Worksheets.Item("loads")

$wsLoads.Range("A1").Value = "name"
$wsLoads.Range("B1").Value = "v_nom_kv"
$wsLoads.Range("C1").Value = "s_base_mva"
$wsLoads.Range("D1").Value = "v_nom_pu"
$wsLoads.Range("E1").Value = "p_nom_mw"
$wsLoads.Range("F1").Value = "q_nom_mvar"
$wsLoads.Range("G1").Value = "bus_idx"
$wsLoads.Range("H1").Value = "g_shunt_pu"
$wsLoads.Range("I1").Value = "b_shunt_pu"

$wsLoads.Range("A2").Value = "Load 1"
$wsLoads.Range("B2").Value = 22
$wsLoads.Range("C2").Value = 100
$wsLoads.Range("D2").Value = 1
$wsLoads.Range("E2").Value = 10
$wsLoads.Range("F2").Value = 10
$wsLoads.Range("G2").Value = 2
$wsLoads.Range("H2").Value = 0
$wsLoads.Range("I2").Value = 0

[void]$wsLoads.Range("J2").Select()

# --- "trafos" sheet (sheet5): extend header with new tap-changer columns ---
$wsTrafos = $wb.Worksheets.Item("trafos")

$wsTrafos.Range("I1").Value = "idx_hv"
$wsTrafos.Range("J1").Value = "idx_lv"
$wsTrafos.Range("K1").Value = "tap_pos"
$wsTrafos.Range("L1").Value = "tap_change"
$wsTrafos.Range("M1").Value = "tap_min"
$wsTrafos.Range("N1").Value = "tap_max"

# trafos becomes the active sheet/tab
[void]$wsTrafos.Activate()
[void]$wsTrafos.Range("A1:N1").Select()
